$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.508.89'
$ws.Range('E2').Value = '  -2.57%  '
$ws.Range('D3').Value = '2.490.03'
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').Value = '315.27'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').Value = '94.01'
$ws.Range('E6').Value = '  -5.09%  '
$ws.Range('D7').Value = '0.549'
$ws.Range('E7').Value = '  -2.32%  '
$ws.Range('E8').Value = '  +0.25%  '
$ws.Range('D9').Value = '0.498'
$ws.Range('E9').Value = '  -3.45%  '
$ws.Range('D10').Value = '33.58'
$ws.Range('E10').Value = '  -4.44%  '
$ws.Range('D11').Value = '0.0785'
$ws.Range('E11').Value = '  -1.82%  '
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('D13').Value = '2.872.55'
$ws.Range('E13').Value = '  -1.14%  '
$ws.Range('D14').Value = '6.93'
$ws.Range('E14').Value = '  -3.64%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '2.566.31'
$ws.Range('E15').Value = '  +3.74%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = '15.52'
$ws.Range('E16').Value = '  +1.80%  '
$ws.Range('E17').Value = '  -1.87%  '
$ws.Range('D18').Value = '41.470.94'
$ws.Range('E18').Value = '  -2.60%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0937'
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '6.37'
$ws.Range('E20').Value = '  -3.38%  '
$ws.Range('D21').Value = '70.34'
$ws.Range('E21').Value = '  +2.19%  '
$ws.Range('D22').Value = '11.25'
$ws.Range('E22').Value = '  -7.02%  '
$ws.Range('D23').Value = '236.66'
$ws.Range('E23').Value = '  -1.73%  '
$ws.Range('D24').Value = '2.78'
$ws.Range('E24').Value = '  -2.62%  '
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('E26').Value = '  -4.52%  '
$ws.Range('D27').Value = '24.28'
$ws.Range('E27').Value = '  -4.28%  '
$ws.Range('E28').Value = '  +0.52%  '
$ws.Range('D29').Value = '9.87'
$ws.Range('E29').Value = '  -1.33%  '
$ws.Range('D30').Value = '37.09'
$ws.Range('D31').Value = '154.38'
$ws.Range('E31').Value = '  -1.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.50'
$ws.Range('E32').Value = '  -5.88%  '
$ws.Range('E33').Value = '  -2.60%  '
$ws.Range('E34').Value = '  -3.67%  '
$ws.Range('D35').Value = '17.93'
$ws.Range('E35').Value = '  +1.77%  '
$ws.Range('E36').Value = '  -2.48%  '
$ws.Range('D37').Value = '2.43'
$ws.Range('E37').Value = '  -9.92%  '
$ws.Range('E38').Value = '  -4.02%  '
$ws.Range('E39').Value = '  -2.67%  '
$ws.Range('E40').Value = '  -6.01%  '
$ws.Range('D41').Value = '4.13'
$ws.Range('E41').Value = '  -1.09%  '
$ws.Range('E42').Value = '  +0.45%  '
$ws.Range('D43').Value = '19.84'
$ws.Range('E43').Value = '  -7.91%  '
$ws.Range('D44').Value = '1.988.77'
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('E45').Value = '  -3.09%  '
$ws.Range('D46').Value = '3.01'
$ws.Range('E46').Value = '  -6.40%  '
$ws.Range('E47').Value = '  -1.03%  '
$ws.Range('D48').Value = '2.733.68'
$ws.Range('E48').Value = '  -0.81%  '
$ws.Range('D49').Value = '69.44'
$ws.Range('E49').Value = '  -2.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '97.00'
$ws.Range('E50').Value = '  -2.75%  '
$ws.Range('E51').Value = '  -5.03%  '
